$d = $word.ActiveDocument

# --- Step 1: remove the existing "_GoBack" bookmark from its old location
# (an otherwise-empty paragraph right after the BIZAGI hyperlink).
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# --- Step 2: rewrite the title paragraph.
# Split the single run into two runs ("...APOIO À " / "COMPLIANCE AND QUALITY
# ASSURANCE") and re-create the "_GoBack" bookmark right after the new text,
# matching the structure Word produces when the title is retyped.
$titlePara = $d.Paragraphs.Item(1).Range
$titleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:r><w:t xml:space="preserve">LISTA DE SOFTWARES DE APOIO À </w:t></w:r>' +
            '<w:r><w:t>COMPLIANCE AND QUALITY ASSURANCE</w:t></w:r>' +
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
            '<w:bookmarkEnd w:id="0"/>' +
            '</w:p>'
$titlePara.InsertXML($titleXml)
